$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3715095283235428
$ws.Cells.Item(2, 3).Value = 0.02941557135680739
$ws.Cells.Item(2, 4).Value = 0.6285659188583423
$ws.Cells.Item(2, 5).Value = 0.2468771440556594
$ws.Cells.Item(2, 7).Value = 0.8779297582206596
$ws.Cells.Item(2, 8).Value = 0.9421776594822688
$ws.Cells.Item(2, 10).Value = 0.121412825833346
$ws.Cells.Item(2, 11).Value = 0.3239916845480195
$ws.Cells.Item(2, 13).Value = 0.2849362235376702
$ws.Cells.Item(2, 14).Value = 1.878773587644797
$ws.Cells.Item(2, 15).Value = 3.671517319374317

$ws.Cells.Item(3, 2).Value = 0.3394262375623214
$ws.Cells.Item(3, 3).Value = 0.02591188334564265
$ws.Cells.Item(3, 4).Value = 0.6229965012113752
$ws.Cells.Item(3, 5).Value = 0.245470005254127
$ws.Cells.Item(3, 7).Value = 0.8805417676847966
$ws.Cells.Item(3, 8).Value = 0.9472367143892626
$ws.Cells.Item(3, 10).Value = 0.1212489028316561
$ws.Cells.Item(3, 11).Value = 0.2910512733943733
$ws.Cells.Item(3, 13).Value = 0.2738213232506226
$ws.Cells.Item(3, 14).Value = 1.896954457759766
$ws.Cells.Item(3, 15).Value = 3.687370667653767

$ws.Cells.Item(4, 2).Value = 0.3198169061441547
$ws.Cells.Item(4, 3).Value = 0.02374964885822806
$ws.Cells.Item(4, 4).Value = 0.6198660070085111
$ws.Cells.Item(4, 5).Value = 0.2447172885537618
$ws.Cells.Item(4, 7).Value = 0.8826113016939914
$ws.Cells.Item(4, 8).Value = 0.9506909768243972
$ws.Cells.Item(4, 10).Value = 0.1212005379292656
$ws.Cells.Item(4, 11).Value = 0.2708641822231499
$ws.Cells.Item(4, 13).Value = 0.2671198289236756
$ws.Cells.Item(4, 14).Value = 1.90868621691865
$ws.Cells.Item(4, 15).Value = 3.698810760536844

$ws.Cells.Item(5, 2).Value = 0.3118490420696958
$ws.Cells.Item(5, 3).Value = 0.022865814405435
$ws.Cells.Item(5, 4).Value = 0.6186631850301154
$ws.Cells.Item(5, 5).Value = 0.244438584074846
$ws.Cells.Item(5, 7).Value = 0.8835717987079121
$ws.Cells.Item(5, 8).Value = 0.9521862194730915
$ws.Cells.Item(5, 10).Value = 0.1211939933223967
$ws.Cells.Item(5, 11).Value = 0.2626479240160364
$ws.Cells.Item(5, 13).Value = 0.264420053808422
$ws.Cells.Item(5, 14).Value = 1.913610017459161
$ws.Cells.Item(5, 15).Value = 3.703901959643574

$ws.Cells.Item(6, 2).Value = 0.3105273935356934
$ws.Cells.Item(6, 3).Value = 0.02271889234300062
$ws.Cells.Item(6, 4).Value = 0.6184678641365196
$ws.Cells.Item(6, 5).Value = 0.2443940001455687
$ws.Cells.Item(6, 7).Value = 0.8837383646321655
$ws.Cells.Item(6, 8).Value = 0.9524397972992205
$ws.Cells.Item(6, 10).Value = 0.1211937022023406
$ws.Cells.Item(6, 11).Value = 0.2612842473579207
$ws.Cells.Item(6, 13).Value = 0.2639736442918803
$ws.Cells.Item(6, 14).Value = 1.914436249262504
$ws.Cells.Item(6, 15).Value = 3.704773283884577

$ws.Cells.Item(7, 2).Value = 0.3197093545823577
$ws.Cells.Item(7, 3).Value = 0.02373774005117468
$ws.Cells.Item(7, 4).Value = 0.6198494900160085
$ws.Cells.Item(7, 5).Value = 0.244713416269235
$ws.Cells.Item(7, 7).Value = 0.8826237809388644
$ws.Cells.Item(7, 8).Value = 0.9507107873584175
$ws.Cells.Item(7, 10).Value = 0.1212003963413011
$ws.Cells.Item(7, 11).Value = 0.2707533331830945
$ws.Cells.Item(7, 13).Value = 0.267083292490824
$ws.Cells.Item(7, 14).Value = 1.908752041923715
$ws.Cells.Item(7, 15).Value = 3.698877683790769

$ws.Cells.Item(8, 2).Value = 0.360428821778072
$ws.Cells.Item(8, 3).Value = 0.02820980175256693
$ws.Cells.Item(8, 4).Value = 0.6265856695733163
$ws.Cells.Item(8, 5).Value = 0.2463688990811725
$ws.Cells.Item(8, 7).Value = 0.87873372826067
$ws.Cells.Item(8, 8).Value = 0.9438498572302905
$ws.Cells.Item(8, 10).Value = 0.1213454634850102
$ws.Cells.Item(8, 11).Value = 0.312626124723181
$ws.Cells.Item(8, 13).Value = 0.2810783630069835
$ws.Cells.Item(8, 14).Value = 1.884924349793269
$ws.Cells.Item(8, 15).Value = 3.676629611429149

$ws.Cells.Item(9, 2).Value = 0.4409760891085739
$ws.Cells.Item(9, 3).Value = 0.03689083427038042
$ws.Cells.Item(9, 4).Value = 0.6420834288017829
$ws.Cells.Item(9, 5).Value = 0.2504962757180067
$ws.Cells.Item(9, 7).Value = 0.8748004184605378
$ws.Cells.Item(9, 8).Value = 0.9331525213841871
$ws.Cells.Item(9, 10).Value = 0.1220441856949961
$ws.Cells.Item(9, 11).Value = 0.395026814746501
$ws.Cells.Item(9, 13).Value = 0.3094929338135657
$ws.Cells.Item(9, 14).Value = 1.842708254159396
$ws.Cells.Item(9, 15).Value = 3.646528802708218

$ws.Cells.Item(10, 2).Value = 0.5005617745589745
$ws.Cells.Item(10, 3).Value = 0.04321298795188966
$ws.Cells.Item(10, 4).Value = 0.6548582539437291
$ws.Cells.Item(10, 5).Value = 0.2540637895852598
$ws.Cells.Item(10, 7).Value = 0.8741638404436713
$ws.Cells.Item(10, 8).Value = 0.9269686298590614
$ws.Cells.Item(10, 10).Value = 0.1228094688008952
$ws.Cells.Item(10, 11).Value = 0.4557256736563886
$ws.Cells.Item(10, 13).Value = 0.3309543387618348
$ws.Cells.Item(10, 14).Value = 1.814437825516063
$ws.Cells.Item(10, 15).Value = 3.6326512431813

$ws.Cells.Item(11, 2).Value = 0.5277540950734476
$ws.Cells.Item(11, 3).Value = 0.04607664718434989
$ws.Cells.Item(11, 4).Value = 0.660970167002688
$ws.Cells.Item(11, 5).Value = 0.2558025712936427
$ws.Cells.Item(11, 7).Value = 0.8743637284049015
$ws.Cells.Item(11, 8).Value = 0.9245181475594535
$ws.Cells.Item(11, 10).Value = 0.1232122069944168
$ws.Cells.Item(11, 11).Value = 0.4833705110629865
$ws.Cells.Item(11, 13).Value = 0.3408435689167035
$ws.Cells.Item(11, 14).Value = 1.802172311105992
$ws.Cells.Item(11, 15).Value = 3.628124977542711

$ws.Cells.Item(12, 2).Value = 0.5380631439409171
$ws.Cells.Item(12, 3).Value = 0.04715922691782737
$ws.Cells.Item(12, 4).Value = 0.6633276689531158
$ws.Cells.Item(12, 5).Value = 0.2564776249779683
$ws.Cells.Item(12, 7).Value = 0.8745098035458057
$ws.Cells.Item(12, 8).Value = 0.923642265651651
$ws.Cells.Item(12, 10).Value = 0.1233725523191822
$ws.Cells.Item(12, 11).Value = 0.4938431757239243
$ws.Cells.Item(12, 13).Value = 0.3446063742570757
$ws.Cells.Item(12, 14).Value = 1.797613197584601
$ws.Cells.Item(12, 15).Value = 3.626667738718453

$ws.Cells.Item(13, 2).Value = 0.5358423846892038
$ws.Cells.Item(13, 3).Value = 0.04692615607108053
$ws.Cells.Item(13, 4).Value = 0.6628180266411903
$ws.Cells.Item(13, 5).Value = 0.2563315019671037
$ws.Cells.Item(13, 7).Value = 0.8744752134534366
$ws.Cells.Item(13, 8).Value = 0.9238285882488952
$ws.Cells.Item(13, 10).Value = 0.1233376707708231
$ws.Cells.Item(13, 11).Value = 0.4915875231674818
$ws.Cells.Item(13, 13).Value = 0.3437951907972021
$ws.Cells.Item(13, 14).Value = 1.798591277933653
$ws.Cells.Item(13, 15).Value = 3.626970164360984

$ws.Cells.Item(14, 2).Value = 0.5286019911301025
$ws.Cells.Item(14, 3).Value = 0.04616574857398348
$ws.Cells.Item(14, 4).Value = 0.6611632582334437
$ws.Cells.Item(14, 5).Value = 0.2558577756781659
$ws.Cells.Item(14, 7).Value = 0.8743743355570359
$ws.Cells.Item(14, 8).Value = 0.924445045267305
$ws.Cells.Item(14, 10).Value = 0.1232252417208386
$ws.Cells.Item(14, 11).Value = 0.484232023037066
$ws.Cells.Item(14, 13).Value = 0.3411527781938872
$ws.Cells.Item(14, 14).Value = 1.801795514934888
$ws.Cells.Item(14, 15).Value = 3.627999944732323

$ws.Cells.Item(15, 2).Value = 0.5241685759343113
$ws.Cells.Item(15, 3).Value = 0.04569973763584301
$ws.Cells.Item(15, 4).Value = 0.6601552665362931
$ws.Cells.Item(15, 5).Value = 0.2555697668457881
$ws.Cells.Item(15, 7).Value = 0.8743217106286636
$ws.Cells.Item(15, 8).Value = 0.9248294207359464
$ws.Cells.Item(15, 10).Value = 0.1231573958506331
$ws.Cells.Item(15, 11).Value = 0.4797270935611948
$ws.Cells.Item(15, 13).Value = 0.3395365583370165
$ws.Cells.Item(15, 14).Value = 1.803769347445872
$ws.Cells.Item(15, 15).Value = 3.628664147707354

$ws.Cells.Item(16, 2).Value = 0.4987863876047811
$ws.Cells.Item(16, 3).Value = 0.04302558853527216
$ws.Cells.Item(16, 4).Value = 0.6544648599825962
$ws.Cells.Item(16, 5).Value = 0.2539524839307319
$ws.Cells.Item(16, 7).Value = 0.8741606264146782
$ws.Cells.Item(16, 8).Value = 0.92713606433243
$ws.Cells.Item(16, 10).Value = 0.1227842464373765
$ws.Cells.Item(16, 11).Value = 0.4539196299762693
$ws.Cells.Item(16, 13).Value = 0.3303105795102894
$ws.Cells.Item(16, 14).Value = 1.81525137826051
$ws.Cells.Item(16, 15).Value = 3.632982984145656

$ws.Cells.Item(17, 2).Value = 0.4832370130318679
$ws.Cells.Item(17, 3).Value = 0.04138189057425734
$ws.Cells.Item(17, 4).Value = 0.6510508450184318
$ws.Cells.Item(17, 5).Value = 0.2529899816883656
$ws.Cells.Item(17, 7).Value = 0.8741871693305399
$ws.Cells.Item(17, 8).Value = 0.9286439298882527
$ws.Cells.Item(17, 10).Value = 0.1225693066454667
$ws.Cells.Item(17, 11).Value = 0.4380955679239094
$ws.Cells.Item(17, 13).Value = 0.3246829559198616
$ws.Cells.Item(17, 14).Value = 1.822447608763154
$ws.Cells.Item(17, 15).Value = 3.636089950170657

$ws.Cells.Item(18, 2).Value = 0.4743015969507383
$ws.Cells.Item(18, 3).Value = 0.04043532258610583
$ws.Cells.Item(18, 4).Value = 0.6491154962144776
$ws.Cells.Item(18, 5).Value = 0.2524472881797664
$ws.Cells.Item(18, 7).Value = 0.8742485076133164
$ws.Cells.Item(18, 8).Value = 0.9295453486821117
$ws.Cells.Item(18, 10).Value = 0.1224508197549738
$ws.Cells.Item(18, 11).Value = 0.4289970879911209
$ws.Cells.Item(18, 13).Value = 0.3214579982306276
$ws.Cells.Item(18, 14).Value = 1.82664267775804
$ws.Cells.Item(18, 15).Value = 3.638045177974362

$ws.Cells.Item(19, 2).Value = 0.4712776384530741
$ws.Cells.Item(19, 3).Value = 0.04011463397212367
$ws.Cells.Item(19, 4).Value = 0.6484650864858281
$ws.Cells.Item(19, 5).Value = 0.2522654172591388
$ws.Cells.Item(19, 7).Value = 0.8742771886242906
$ws.Cells.Item(19, 8).Value = 0.9298564186320135
$ws.Cells.Item(19, 10).Value = 0.122411585532177
$ws.Cells.Item(19, 11).Value = 0.4259170478253225
$ws.Cells.Item(19, 13).Value = 0.3203681332559682
$ws.Cells.Item(19, 14).Value = 1.82807267100704
$ws.Cells.Item(19, 15).Value = 3.638736075206452

$ws.Cells.Item(20, 2).Value = 0.4848914284744126
$ws.Cells.Item(20, 3).Value = 0.04155698510564321
$ws.Cells.Item(20, 4).Value = 0.6514113447882153
$ws.Cells.Item(20, 5).Value = 0.2530913126853918
$ws.Cells.Item(20, 7).Value = 0.874179575587263
$ws.Cells.Item(20, 8).Value = 0.9284798828160632
$ws.Cells.Item(20, 10).Value = 0.1225916554110782
$ws.Cells.Item(20, 11).Value = 0.4397797491062079
$ws.Cells.Item(20, 13).Value = 0.3252807956181556
$ws.Cells.Item(20, 14).Value = 1.821675763011221
$ws.Cells.Item(20, 15).Value = 3.635741803214103

$ws.Cells.Item(21, 2).Value = 0.5307283512997003
$ws.Cells.Item(21, 3).Value = 0.04638914867460642
$ws.Cells.Item(21, 4).Value = 0.6616481367348399
$ws.Cells.Item(21, 5).Value = 0.2559964700649573
$ws.Cells.Item(21, 7).Value = 0.8744020557240617
$ws.Cells.Item(21, 8).Value = 0.9242625645618716
$ws.Cells.Item(21, 10).Value = 0.1232580522836173
$ws.Cells.Item(21, 11).Value = 0.4863924028929887
$ws.Cells.Item(21, 13).Value = 0.3419284327445737
$ws.Cells.Item(21, 14).Value = 1.800852030222049
$ws.Cells.Item(21, 15).Value = 3.627690506214236

$ws.Cells.Item(22, 2).Value = 0.5607545879947509
$ws.Cells.Item(22, 3).Value = 0.04953657920719934
$ws.Cells.Item(22, 4).Value = 0.6685893393439812
$ws.Cells.Item(22, 5).Value = 0.2579919779529547
$ws.Cells.Item(22, 7).Value = 0.874957697163822
$ws.Cells.Item(22, 8).Value = 0.9218097274013815
$ws.Cells.Item(22, 10).Value = 0.1237392547242848
$ws.Cells.Item(22, 11).Value = 0.5168804972618943
$ws.Cells.Item(22, 13).Value = 0.352913268217641
$ws.Cells.Item(22, 14).Value = 1.787741304894492
$ws.Cells.Item(22, 15).Value = 3.62392503341556

$ws.Cells.Item(23, 2).Value = 0.5447228678466729
$ws.Cells.Item(23, 3).Value = 0.04785773030256735
$ws.Cells.Item(23, 4).Value = 0.6648617908069241
$ws.Cells.Item(23, 5).Value = 0.2569180959659292
$ws.Cells.Item(23, 7).Value = 0.8746236052860894
$ws.Cells.Item(23, 8).Value = 0.923091115409008
$ws.Cells.Item(23, 10).Value = 0.1234782538571082
$ws.Cells.Item(23, 11).Value = 0.5006064027297441
$ws.Cells.Item(23, 13).Value = 0.347040945385956
$ws.Cells.Item(23, 14).Value = 1.794693095199436
$ws.Cells.Item(23, 15).Value = 3.625797860786662

$ws.Cells.Item(24, 2).Value = 0.4841434543598098
$ws.Cells.Item(24, 3).Value = 0.04147782979781311
$ws.Cells.Item(24, 4).Value = 0.6512482774391515
$ws.Cells.Item(24, 5).Value = 0.2530454677276666
$ws.Cells.Item(24, 7).Value = 0.8741828651868389
$ws.Cells.Item(24, 8).Value = 0.9285539409410717
$ws.Cells.Item(24, 10).Value = 0.1225815356946427
$ws.Cells.Item(24, 11).Value = 0.4390183339185683
$ws.Cells.Item(24, 13).Value = 0.3250104797627813
$ws.Cells.Item(24, 14).Value = 1.822024534244469
$ws.Cells.Item(24, 15).Value = 3.635898674039652

$ws.Cells.Item(25, 2).Value = 0.4191131085464121
$ws.Cells.Item(25, 3).Value = 0.03455205817373042
$ws.Cells.Item(25, 4).Value = 0.6376465625555738
$ws.Cells.Item(25, 5).Value = 0.2492855961220108
$ws.Cells.Item(25, 7).Value = 0.8754688365147985
$ws.Cells.Item(25, 8).Value = 0.9357518316621309
$ws.Cells.Item(25, 10).Value = 0.1218108773346387
$ws.Cells.Item(25, 11).Value = 0.3727061050395264
$ws.Cells.Item(25, 13).Value = 0.3309543387618348
$ws.Cells.Item(25, 14).Value = 1.853646393820312
$ws.Cells.Item(25, 15).Value = 3.653224609546868
